# Update the metrics table (columns B:Q, rows 2:26) with new computed values.
# All rows share the same values per column, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B" = 0.6383931775788736
    "C" = -23.10748969621655
    "D" = -1.13872184183592
    "E" = 0.4959262137176474
    "F" = 0.1584429363339767
    "G" = 0.2146654303905574
    "H" = 14.31124727300474
    "I" = 0.1890308038581712
    "J" = 0.1706313204555426
    "K" = 0.1798310621568569
    "L" = 0.2718993657310753
    "M" = 0.4633200086231518
    "N" = -0.08482046726337922
    "O" = 0.4830445005948765
    "P" = 35.07734920122535
    "Q" = 54.57936239911656
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
